# Form the consolidated report: fill in / correct the "Absent" column (H)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H12").Value = 0
